# Applies the "added indent to match style" edit to the styles template docx.
# Targets: numbering.xml abstractNum#9 tmpl value, and the Heading1-4 / HeadingNChar /
# Indent1-4 style definitions inside styles.xml (rsid bump, tab-stop clears, hanging
# indents on Heading2/3, new <w:lang w:eastAsia="zh-CN"/> runs, and revised Indent2-4 left indents).

$d = $word.ActiveDocument

# Word exposes the full package as flattened WordprocessingML via Document.WordOpenXML;
# word/styles.xml and word/numbering.xml both round-trip through it, so editing the
# text there and writing it back mutates those two parts precisely.
$xml = $d.WordOpenXML

# --- numbering.xml: abstractNum 9 tmpl id ---
$oldTmpl = '<w:tmpl w:val="19E23370"/>'
$newTmpl = '<w:tmpl w:val="B42C7DE0"/>'
if ($xml.IndexOf($oldTmpl) -lt 0) { throw "tmpl anchor not found" }
$xml = $xml.Replace($oldTmpl, $newTmpl)

# --- styles.xml: Heading1 ---
$oldHeading1 = '<w:style w:type="paragraph" w:styleId="Heading1"><w:name w:val="heading 1"/><w:basedOn w:val="Normal"/><w:next w:val="Normal"/><w:link w:val="Heading1Char"/><w:uiPriority w:val="9"/><w:qFormat/><w:rsid w:val="007A00B2"/><w:pPr><w:widowControl w:val="0"/><w:numPr><w:numId w:val="38"/></w:numPr><w:spacing w:before="240" w:after="0" w:line="240" w:lineRule="auto"/><w:outlineLvl w:val="0"/></w:pPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsia="SimSun" w:hAnsiTheme="majorHAnsi" w:cs="Arial"/><w:b/><w:bCs/><w:kern w:val="32"/><w:sz w:val="28"/><w:szCs w:val="32"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:style>'
$newHeading1 = '<w:style w:type="paragraph" w:styleId="Heading1"><w:name w:val="heading 1"/><w:basedOn w:val="Normal"/><w:next w:val="Normal"/><w:link w:val="Heading1Char"/><w:uiPriority w:val="9"/><w:qFormat/><w:rsid w:val="00DF71D1"/><w:pPr><w:widowControl w:val="0"/><w:numPr><w:numId w:val="38"/></w:numPr><w:tabs><w:tab w:val="clear" w:pos="360"/></w:tabs><w:spacing w:before="240" w:after="0" w:line="240" w:lineRule="auto"/><w:outlineLvl w:val="0"/></w:pPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsia="SimSun" w:hAnsiTheme="majorHAnsi" w:cs="Arial"/><w:b/><w:bCs/><w:kern w:val="32"/><w:sz w:val="28"/><w:szCs w:val="32"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:style>'
if ($xml.IndexOf($oldHeading1) -lt 0) { throw "Heading1 anchor not found" }
$xml = $xml.Replace($oldHeading1, $newHeading1)

# --- styles.xml: Heading2 ---
$oldHeading2 = '<w:style w:type="paragraph" w:styleId="Heading2"><w:name w:val="heading 2"/><w:basedOn w:val="Normal"/><w:next w:val="Normal"/><w:link w:val="Heading2Char"/><w:uiPriority w:val="9"/><w:qFormat/><w:rsid w:val="007A00B2"/><w:pPr><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="38"/></w:numPr><w:spacing w:before="60" w:line="240" w:lineRule="auto"/><w:outlineLvl w:val="1"/></w:pPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsia="SimSun" w:hAnsiTheme="majorHAnsi" w:cs="Arial"/><w:b/><w:bCs/><w:iCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:style>'
$newHeading2 = '<w:style w:type="paragraph" w:styleId="Heading2"><w:name w:val="heading 2"/><w:basedOn w:val="Normal"/><w:next w:val="Normal"/><w:link w:val="Heading2Char"/><w:uiPriority w:val="9"/><w:qFormat/><w:rsid w:val="00DF71D1"/><w:pPr><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="38"/></w:numPr><w:tabs><w:tab w:val="clear" w:pos="1080"/></w:tabs><w:spacing w:before="60" w:line="240" w:lineRule="auto"/><w:ind w:left="864" w:hanging="504"/><w:outlineLvl w:val="1"/></w:pPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsia="SimSun" w:hAnsiTheme="majorHAnsi" w:cs="Arial"/><w:b/><w:bCs/><w:iCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:style>'
if ($xml.IndexOf($oldHeading2) -lt 0) { throw "Heading2 anchor not found" }
$xml = $xml.Replace($oldHeading2, $newHeading2)

# --- styles.xml: Heading3 ---
$oldHeading3 = '<w:style w:type="paragraph" w:styleId="Heading3"><w:name w:val="heading 3"/><w:basedOn w:val="Normal"/><w:next w:val="Normal"/><w:link w:val="Heading3Char"/><w:uiPriority w:val="9"/><w:qFormat/><w:rsid w:val="007A00B2"/><w:pPr><w:numPr><w:ilvl w:val="2"/><w:numId w:val="38"/></w:numPr><w:spacing w:before="60" w:line="240" w:lineRule="auto"/><w:outlineLvl w:val="2"/></w:pPr><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="SimSun" w:hAnsi="Cambria" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:style>'
$newHeading3 = '<w:style w:type="paragraph" w:styleId="Heading3"><w:name w:val="heading 3"/><w:basedOn w:val="Normal"/><w:next w:val="Normal"/><w:link w:val="Heading3Char"/><w:uiPriority w:val="9"/><w:qFormat/><w:rsid w:val="00DF71D1"/><w:pPr><w:numPr><w:ilvl w:val="2"/><w:numId w:val="38"/></w:numPr><w:tabs><w:tab w:val="clear" w:pos="1800"/></w:tabs><w:spacing w:before="60" w:line="240" w:lineRule="auto"/><w:ind w:left="1440" w:hanging="720"/><w:outlineLvl w:val="2"/></w:pPr><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="SimSun" w:hAnsi="Cambria" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:style>'
if ($xml.IndexOf($oldHeading3) -lt 0) { throw "Heading3 anchor not found" }
$xml = $xml.Replace($oldHeading3, $newHeading3)

# --- styles.xml: Heading4 ---
$oldHeading4 = '<w:style w:type="paragraph" w:styleId="Heading4"><w:name w:val="heading 4"/><w:basedOn w:val="Normal"/><w:next w:val="Normal"/><w:link w:val="Heading4Char"/><w:uiPriority w:val="9"/><w:unhideWhenUsed/><w:qFormat/><w:rsid w:val="007A00B2"/><w:pPr><w:keepNext/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="38"/></w:numPr><w:spacing w:before="120"/><w:outlineLvl w:val="3"/></w:pPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="minorBidi"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr></w:style>'
$newHeading4 = '<w:style w:type="paragraph" w:styleId="Heading4"><w:name w:val="heading 4"/><w:basedOn w:val="Normal"/><w:next w:val="Normal"/><w:link w:val="Heading4Char"/><w:uiPriority w:val="9"/><w:unhideWhenUsed/><w:qFormat/><w:rsid w:val="00DF71D1"/><w:pPr><w:keepNext/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="38"/></w:numPr><w:tabs><w:tab w:val="clear" w:pos="1800"/></w:tabs><w:spacing w:before="120"/><w:outlineLvl w:val="3"/></w:pPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="minorBidi"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="28"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:style>'
if ($xml.IndexOf($oldHeading4) -lt 0) { throw "Heading4 anchor not found" }
$xml = $xml.Replace($oldHeading4, $newHeading4)

# --- styles.xml: Heading1Char ---
$oldHeading1Char = '<w:style w:type="character" w:customStyle="1" w:styleId="Heading1Char"><w:name w:val="Heading 1 Char"/><w:basedOn w:val="DefaultParagraphFont"/><w:link w:val="Heading1"/><w:uiPriority w:val="9"/><w:rsid w:val="007A00B2"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsia="SimSun" w:hAnsiTheme="majorHAnsi" w:cs="Arial"/><w:b/><w:bCs/><w:kern w:val="32"/><w:sz w:val="28"/><w:szCs w:val="32"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:style>'
$newHeading1Char = '<w:style w:type="character" w:customStyle="1" w:styleId="Heading1Char"><w:name w:val="Heading 1 Char"/><w:basedOn w:val="DefaultParagraphFont"/><w:link w:val="Heading1"/><w:uiPriority w:val="9"/><w:rsid w:val="00DF71D1"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsia="SimSun" w:hAnsiTheme="majorHAnsi" w:cs="Arial"/><w:b/><w:bCs/><w:kern w:val="32"/><w:sz w:val="28"/><w:szCs w:val="32"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:style>'
if ($xml.IndexOf($oldHeading1Char) -lt 0) { throw "Heading1Char anchor not found" }
$xml = $xml.Replace($oldHeading1Char, $newHeading1Char)

# --- styles.xml: Heading2Char ---
$oldHeading2Char = '<w:style w:type="character" w:customStyle="1" w:styleId="Heading2Char"><w:name w:val="Heading 2 Char"/><w:basedOn w:val="DefaultParagraphFont"/><w:link w:val="Heading2"/><w:uiPriority w:val="9"/><w:rsid w:val="007A00B2"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsia="SimSun" w:hAnsiTheme="majorHAnsi" w:cs="Arial"/><w:b/><w:bCs/><w:iCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:style>'
$newHeading2Char = '<w:style w:type="character" w:customStyle="1" w:styleId="Heading2Char"><w:name w:val="Heading 2 Char"/><w:basedOn w:val="DefaultParagraphFont"/><w:link w:val="Heading2"/><w:uiPriority w:val="9"/><w:rsid w:val="00DF71D1"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsia="SimSun" w:hAnsiTheme="majorHAnsi" w:cs="Arial"/><w:b/><w:bCs/><w:iCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:style>'
if ($xml.IndexOf($oldHeading2Char) -lt 0) { throw "Heading2Char anchor not found" }
$xml = $xml.Replace($oldHeading2Char, $newHeading2Char)

# --- styles.xml: Heading3Char ---
$oldHeading3Char = '<w:style w:type="character" w:customStyle="1" w:styleId="Heading3Char"><w:name w:val="Heading 3 Char"/><w:basedOn w:val="DefaultParagraphFont"/><w:link w:val="Heading3"/><w:uiPriority w:val="9"/><w:rsid w:val="007A00B2"/><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="SimSun" w:hAnsi="Cambria" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:style>'
$newHeading3Char = '<w:style w:type="character" w:customStyle="1" w:styleId="Heading3Char"><w:name w:val="Heading 3 Char"/><w:basedOn w:val="DefaultParagraphFont"/><w:link w:val="Heading3"/><w:uiPriority w:val="9"/><w:rsid w:val="00DF71D1"/><w:rPr><w:rFonts w:ascii="Cambria" w:eastAsia="SimSun" w:hAnsi="Cambria" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:style>'
if ($xml.IndexOf($oldHeading3Char) -lt 0) { throw "Heading3Char anchor not found" }
$xml = $xml.Replace($oldHeading3Char, $newHeading3Char)

# --- styles.xml: Heading4Char ---
$oldHeading4Char = '<w:style w:type="character" w:customStyle="1" w:styleId="Heading4Char"><w:name w:val="Heading 4 Char"/><w:basedOn w:val="DefaultParagraphFont"/><w:link w:val="Heading4"/><w:uiPriority w:val="9"/><w:rsid w:val="007A00B2"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="minorBidi"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr></w:style>'
$newHeading4Char = '<w:style w:type="character" w:customStyle="1" w:styleId="Heading4Char"><w:name w:val="Heading 4 Char"/><w:basedOn w:val="DefaultParagraphFont"/><w:link w:val="Heading4"/><w:uiPriority w:val="9"/><w:rsid w:val="00DF71D1"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="minorBidi"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="28"/><w:lang w:eastAsia="zh-CN"/></w:rPr></w:style>'
if ($xml.IndexOf($oldHeading4Char) -lt 0) { throw "Heading4Char anchor not found" }
$xml = $xml.Replace($oldHeading4Char, $newHeading4Char)

# --- styles.xml: Indent1 ---
$oldIndent1 = '<w:style w:type="paragraph" w:customStyle="1" w:styleId="Indent1"><w:name w:val="Indent 1"/><w:basedOn w:val="Normal"/><w:qFormat/><w:rsid w:val="007A00B2"/><w:pPr><w:ind w:left="360"/></w:pPr></w:style>'
$newIndent1 = '<w:style w:type="paragraph" w:customStyle="1" w:styleId="Indent1"><w:name w:val="Indent 1"/><w:basedOn w:val="Normal"/><w:qFormat/><w:rsid w:val="00DF71D1"/><w:pPr><w:ind w:left="360"/></w:pPr><w:rPr><w:lang w:eastAsia="zh-CN"/></w:rPr></w:style>'
if ($xml.IndexOf($oldIndent1) -lt 0) { throw "Indent1 anchor not found" }
$xml = $xml.Replace($oldIndent1, $newIndent1)

# --- styles.xml: Indent2 ---
$oldIndent2 = '<w:style w:type="paragraph" w:customStyle="1" w:styleId="Indent2"><w:name w:val="Indent 2"/><w:basedOn w:val="Normal"/><w:qFormat/><w:rsid w:val="007A00B2"/><w:pPr><w:spacing w:after="120" w:line="240" w:lineRule="auto"/><w:ind w:left="810"/></w:pPr></w:style>'
$newIndent2 = '<w:style w:type="paragraph" w:customStyle="1" w:styleId="Indent2"><w:name w:val="Indent 2"/><w:basedOn w:val="Normal"/><w:qFormat/><w:rsid w:val="00DF71D1"/><w:pPr><w:spacing w:after="120" w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:rPr><w:lang w:eastAsia="zh-CN"/></w:rPr></w:style>'
if ($xml.IndexOf($oldIndent2) -lt 0) { throw "Indent2 anchor not found" }
$xml = $xml.Replace($oldIndent2, $newIndent2)

# --- styles.xml: Indent3 ---
$oldIndent3 = '<w:style w:type="paragraph" w:customStyle="1" w:styleId="Indent3"><w:name w:val="Indent 3"/><w:basedOn w:val="Normal"/><w:qFormat/><w:rsid w:val="007A00B2"/><w:pPr><w:spacing w:after="120" w:line="240" w:lineRule="auto"/><w:ind w:left="1350"/></w:pPr></w:style>'
$newIndent3 = '<w:style w:type="paragraph" w:customStyle="1" w:styleId="Indent3"><w:name w:val="Indent 3"/><w:basedOn w:val="Normal"/><w:qFormat/><w:rsid w:val="00DF71D1"/><w:pPr><w:spacing w:after="120" w:line="240" w:lineRule="auto"/><w:ind w:left="1080"/></w:pPr><w:rPr><w:lang w:eastAsia="zh-CN"/></w:rPr></w:style>'
if ($xml.IndexOf($oldIndent3) -lt 0) { throw "Indent3 anchor not found" }
$xml = $xml.Replace($oldIndent3, $newIndent3)

# --- styles.xml: Indent4 ---
$oldIndent4 = '<w:style w:type="paragraph" w:customStyle="1" w:styleId="Indent4"><w:name w:val="Indent 4"/><w:basedOn w:val="Normal"/><w:qFormat/><w:rsid w:val="007A00B2"/><w:pPr><w:spacing w:after="120" w:line="240" w:lineRule="auto"/><w:ind w:left="1987"/></w:pPr></w:style>'
$newIndent4 = '<w:style w:type="paragraph" w:customStyle="1" w:styleId="Indent4"><w:name w:val="Indent 4"/><w:basedOn w:val="Normal"/><w:qFormat/><w:rsid w:val="00DF71D1"/><w:pPr><w:spacing w:after="120" w:line="240" w:lineRule="auto"/><w:ind w:left="1440"/></w:pPr><w:rPr><w:lang w:eastAsia="zh-CN"/></w:rPr></w:style>'
if ($xml.IndexOf($oldIndent4) -lt 0) { throw "Indent4 anchor not found" }
$xml = $xml.Replace($oldIndent4, $newIndent4)

# Write the patched package XML back; this commits the styles.xml / numbering.xml edits.
$d.WordOpenXML = $xml

Write-Output "style template updated"
